# Applies the scheduled-runner update to Kujata_Profits workbook.
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) holds market/profit
# data in columns H:N (currentAveragePrice, currentAveragePriceNQ,
# currentAveragePriceHQ, LevePriceNQ, LevePriceHQ, LeveProfitNQ, LeveProfitHQ).
# This script refreshes those figures for the rows the scheduled runner
# re-priced.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ALC
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 96
$ws.Cells.Item(96,8).Value  = 600
$ws.Cells.Item(96,9).Value  = 250
$ws.Cells.Item(96,10).Value = 1300
$ws.Cells.Item(96,11).Value = 750
$ws.Cells.Item(96,12).Value = 3900
$ws.Cells.Item(96,13).Value = 623
$ws.Cells.Item(96,14).Value = -6646

# Row 137
$ws.Cells.Item(137,8).Value  = 1309.4482
$ws.Cells.Item(137,10).Value = 1601.2858
$ws.Cells.Item(137,12).Value = 4803.857400000001
$ws.Cells.Item(137,14).Value = -9903.857400000001

# ---------------------------------------------------------------------
# ARM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 34
$ws.Cells.Item(34,8).Value  = 7352
$ws.Cells.Item(34,9).Value  = 2000
$ws.Cells.Item(34,11).Value = 2000
$ws.Cells.Item(34,13).Value = -1729

# Row 132
$ws.Cells.Item(132,8).Value  = 2415.5715
$ws.Cells.Item(132,9).Value  = 2120.9375
$ws.Cells.Item(132,10).Value = 2808.4167
$ws.Cells.Item(132,11).Value = 6362.8125
$ws.Cells.Item(132,12).Value = 8425.250100000001
$ws.Cells.Item(132,13).Value = -3832.8125
$ws.Cells.Item(132,14).Value = -13485.2501

# ---------------------------------------------------------------------
# BSM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 105
$ws.Cells.Item(105,8).Value  = 76925210
$ws.Cells.Item(105,9).Value  = 125001930
$ws.Cells.Item(105,11).Value = 125001930
$ws.Cells.Item(105,13).Value = -125000183

# Row 134
$ws.Cells.Item(134,8).Value  = 7308.05
$ws.Cells.Item(134,9).Value  = 1011.6429
$ws.Cells.Item(134,10).Value = 21999.666
$ws.Cells.Item(134,11).Value = 3034.9287
$ws.Cells.Item(134,12).Value = 65998.99800000001
$ws.Cells.Item(134,13).Value = -499.9287000000004
$ws.Cells.Item(134,14).Value = -71068.99800000001

# ---------------------------------------------------------------------
# CRP
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 16
$ws.Cells.Item(16,8).Value  = 83334350
$ws.Cells.Item(16,9).Value  = 100000990
$ws.Cells.Item(16,10).Value = 1144
$ws.Cells.Item(16,11).Value = 100000990
$ws.Cells.Item(16,12).Value = 1144
$ws.Cells.Item(16,13).Value = -100000703
$ws.Cells.Item(16,14).Value = -1718

# Row 22
$ws.Cells.Item(22,8).Value  = 437.5
$ws.Cells.Item(22,9).Value  = 414
$ws.Cells.Item(22,10).Value = 476.66666
$ws.Cells.Item(22,11).Value = 414
$ws.Cells.Item(22,12).Value = 476.66666
$ws.Cells.Item(22,13).Value = -64
$ws.Cells.Item(22,14).Value = -1176.66666

# Row 113
$ws.Cells.Item(113,8).Value  = 83334350
$ws.Cells.Item(113,9).Value  = 100000990
$ws.Cells.Item(113,10).Value = 1144
$ws.Cells.Item(113,11).Value = 100000990
$ws.Cells.Item(113,12).Value = 1144
$ws.Cells.Item(113,13).Value = -99998820
$ws.Cells.Item(113,14).Value = -5484

# ---------------------------------------------------------------------
# CUL
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 110
$ws.Cells.Item(110,8).Value = 10753.375

# Row 116
$ws.Cells.Item(116,8).Value  = 835208.3
$ws.Cells.Item(116,9).Value  = 2500200
$ws.Cells.Item(116,11).Value = 7500600
$ws.Cells.Item(116,13).Value = -7497158

# ---------------------------------------------------------------------
# GSM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 11
$ws.Cells.Item(11,8).Value  = 7179167
$ws.Cells.Item(11,9).Value  = 7568750
$ws.Cells.Item(11,10).Value = 4062500
$ws.Cells.Item(11,11).Value = 7568750
$ws.Cells.Item(11,12).Value = 4062500
$ws.Cells.Item(11,13).Value = -7568611
$ws.Cells.Item(11,14).Value = -4062778

# ---------------------------------------------------------------------
# LTW
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 38 - prices reset to 0, LeveProfitNQ (M) cleared entirely
$ws.Cells.Item(38,8).Value  = 0
$ws.Cells.Item(38,9).Value  = 0
$ws.Cells.Item(38,11).Value = 0
$ws.Cells.Item(38,13).ClearContents()

# Row 124
$ws.Cells.Item(124,8).Value  = 0
$ws.Cells.Item(124,9).Value  = 0
$ws.Cells.Item(124,10).Value = 0
$ws.Cells.Item(124,11).Value = 0
$ws.Cells.Item(124,12).Value = 0

# Row 125
$ws.Cells.Item(125,8).Value  = 10157.5
$ws.Cells.Item(125,9).Value  = 0
$ws.Cells.Item(125,10).Value = 10157.5
$ws.Cells.Item(125,11).Value = 0
$ws.Cells.Item(125,12).Value = 10157.5
$ws.Cells.Item(125,14).Value = -19997.5

# Row 127
$ws.Cells.Item(127,8).Value  = 33000
$ws.Cells.Item(127,9).Value  = 0
$ws.Cells.Item(127,10).Value = 33000
$ws.Cells.Item(127,11).Value = 0
$ws.Cells.Item(127,12).Value = 33000
$ws.Cells.Item(127,14).Value = -42920

# Row 128
$ws.Cells.Item(128,8).Value  = 90000
$ws.Cells.Item(128,9).Value  = 0
$ws.Cells.Item(128,10).Value = 90000
$ws.Cells.Item(128,11).Value = 0
$ws.Cells.Item(128,12).Value = 90000
$ws.Cells.Item(128,14).Value = -99960

# Row 129
$ws.Cells.Item(129,8).Value  = 0
$ws.Cells.Item(129,9).Value  = 0
$ws.Cells.Item(129,10).Value = 0
$ws.Cells.Item(129,11).Value = 0
$ws.Cells.Item(129,12).Value = 0

# Row 130
$ws.Cells.Item(130,8).Value  = 71552
$ws.Cells.Item(130,9).Value  = 0
$ws.Cells.Item(130,10).Value = 71552
$ws.Cells.Item(130,11).Value = 0
$ws.Cells.Item(130,12).Value = 71552
$ws.Cells.Item(130,14).Value = -81592

# Row 131
$ws.Cells.Item(131,8).Value  = 50000
$ws.Cells.Item(131,9).Value  = 0
$ws.Cells.Item(131,10).Value = 50000
$ws.Cells.Item(131,11).Value = 0
$ws.Cells.Item(131,12).Value = 50000
$ws.Cells.Item(131,14).Value = -60080

# Row 132
$ws.Cells.Item(132,8).Value  = 21989.918
$ws.Cells.Item(132,9).Value  = 1335.5
$ws.Cells.Item(132,10).Value = 49529.145
$ws.Cells.Item(132,11).Value = 4006.5
$ws.Cells.Item(132,12).Value = 148587.435
$ws.Cells.Item(132,13).Value = -1476.5
$ws.Cells.Item(132,14).Value = -153647.435

# Row 133
$ws.Cells.Item(133,8).Value  = 34019.8
$ws.Cells.Item(133,9).Value  = 0
$ws.Cells.Item(133,10).Value = 34019.8
$ws.Cells.Item(133,11).Value = 0
$ws.Cells.Item(133,12).Value = 34019.8
$ws.Cells.Item(133,14).Value = -39079.8

# Row 134
$ws.Cells.Item(134,8).Value  = 0
$ws.Cells.Item(134,9).Value  = 0
$ws.Cells.Item(134,10).Value = 0
$ws.Cells.Item(134,11).Value = 0
$ws.Cells.Item(134,12).Value = 0

# Row 135
$ws.Cells.Item(135,8).Value  = 0
$ws.Cells.Item(135,9).Value  = 0
$ws.Cells.Item(135,10).Value = 0
$ws.Cells.Item(135,11).Value = 0
$ws.Cells.Item(135,12).Value = 0

# Row 136
$ws.Cells.Item(136,8).Value  = 8306.6
$ws.Cells.Item(136,9).Value  = 12266.667
$ws.Cells.Item(136,10).Value = 2366.5
$ws.Cells.Item(136,11).Value = 36800.001
$ws.Cells.Item(136,12).Value = 7099.5
$ws.Cells.Item(136,13).Value = -34250.001
$ws.Cells.Item(136,14).Value = -12199.5

# Row 137
$ws.Cells.Item(137,8).Value  = 33830
$ws.Cells.Item(137,9).Value  = 0
$ws.Cells.Item(137,10).Value = 33830
$ws.Cells.Item(137,11).Value = 0
$ws.Cells.Item(137,12).Value = 33830
$ws.Cells.Item(137,14).Value = -44030

# Row 138
$ws.Cells.Item(138,8).Value  = 39950
$ws.Cells.Item(138,9).Value  = 0
$ws.Cells.Item(138,10).Value = 39950
$ws.Cells.Item(138,11).Value = 0
$ws.Cells.Item(138,12).Value = 39950
$ws.Cells.Item(138,14).Value = -50230

# Row 139
$ws.Cells.Item(139,8).Value  = 0
$ws.Cells.Item(139,9).Value  = 0
$ws.Cells.Item(139,10).Value = 0
$ws.Cells.Item(139,11).Value = 0
$ws.Cells.Item(139,12).Value = 0

# Row 140
$ws.Cells.Item(140,8).Value  = 38067.418
$ws.Cells.Item(140,9).Value  = 0
$ws.Cells.Item(140,10).Value = 38067.418
$ws.Cells.Item(140,11).Value = 0
$ws.Cells.Item(140,12).Value = 38067.418
$ws.Cells.Item(140,14).Value = -48427.418

# Row 141
$ws.Cells.Item(141,8).Value  = 47928.75
$ws.Cells.Item(141,9).Value  = 0
$ws.Cells.Item(141,10).Value = 47928.75
$ws.Cells.Item(141,11).Value = 0
$ws.Cells.Item(141,12).Value = 47928.75
$ws.Cells.Item(141,14).Value = -58288.75

Write-Host "Kujata_Profits update applied."
